$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed API pull no longer includes the last 11 coins (rows 35-45),
# so those rows are removed and the sheet shrinks to A1:J34.
$ws.Range("A35:J45").EntireRow.Delete()

# Columns holding text-like values (ids, percentages, btc price, supply figures)
# must stay text so figures such as "1.00" or "0.000016" keep their exact
# formatting instead of being coerced into numbers.
$ws.Range("A2:A34").NumberFormat = "@"
$ws.Range("B2:B34").NumberFormat = "@"
$ws.Range("C2:C34").NumberFormat = "@"
$ws.Range("E2:G34").NumberFormat = "@"
$ws.Range("I2:I34").NumberFormat = "@"
$ws.Range("J2:J28").NumberFormat = "@"
$ws.Range("J30:J33").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "90"
$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "bitcoin"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "0.35"
$ws.Range("F2").Value = "0.21"
$ws.Range("G2").Value = "1.00"
$ws.Range("H2").Value = 27467658598.48648
$ws.Range("I2").Value = "19690059.00"
$ws.Range("J2").Value = "21000000"

# Row 3
$ws.Range("A3").Value = "80"
$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "ethereum"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "-0.42"
$ws.Range("F3").Value = "-3.36"
$ws.Range("G3").Value = "0.048270"
$ws.Range("H3").Value = 10771076454.20971
$ws.Range("I3").Value = "122375302.00"
$ws.Range("J3").Value = ""

# Row 4
$ws.Range("A4").Value = "518"
$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "tether"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "-0.05"
$ws.Range("F4").Value = "-0.30"
$ws.Range("G4").Value = "0.000016"
$ws.Range("H4").Value = 44835413287.13708
$ws.Range("I4").Value = "100729057191.00"
$ws.Range("J4").Value = ""

# Row 5
$ws.Range("A5").Value = "2710"
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "binance-coin"
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = "0.84"
$ws.Range("F5").Value = "-1.25"
$ws.Range("G5").Value = "0.009344"
$ws.Range("H5").Value = 501162059.1712992
$ws.Range("I5").Value = "166801148.00"
$ws.Range("J5").Value = "200000000"

# Row 6
$ws.Range("A6").Value = "48543"
$ws.Range("B6").Value = "SOL"
$ws.Range("C6").Value = "solana"
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = "1.29"
$ws.Range("F6").Value = "13.73"
$ws.Range("G6").Value = "0.002423"
$ws.Range("H6").Value = 2499304258.724568
$ws.Range("I6").Value = "443015903.00"
$ws.Range("J6").Value = ""

# Row 7
$ws.Range("A7").Value = "33285"
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "usd-coin"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = "0.01"
$ws.Range("F7").Value = "-0.01"
$ws.Range("G7").Value = "0.000016"
$ws.Range("H7").Value = 3437834993.044488
$ws.Range("I7").Value = "29187435574.00"
$ws.Range("J7").Value = ""

# Row 8
$ws.Range("A8").Value = "46971"
$ws.Range("B8").Value = "STETH"
$ws.Range("C8").Value = "staked-ether"
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = "-0.42"
$ws.Range("F8").Value = "-3.40"
$ws.Range("G8").Value = "0.048193"
$ws.Range("H8").Value = 81584646.67647275
$ws.Range("I8").Value = "9220200.00"
$ws.Range("J8").Value = "9901295"

# Row 9
$ws.Range("A9").Value = "58"
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "ripple"
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = "-1.86"
$ws.Range("F9").Value = "5.31"
$ws.Range("G9").Value = "0.000009"
$ws.Range("H9").Value = 1443695100.808488
$ws.Range("I9").Value = "42909539227.00"
$ws.Range("J9").Value = "100000000000"

# Row 10
$ws.Range("A10").Value = "2"
$ws.Range("B10").Value = "DOGE"
$ws.Range("C10").Value = "dogecoin"
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = "-0.34"
$ws.Range("F10").Value = "8.73"
$ws.Range("G10").Value = "0.000002"
$ws.Range("H10").Value = 1428667307.267194
$ws.Range("I10").Value = "144031626384.00"
$ws.Range("J10").Value = ""

# Row 11
$ws.Range("A11").Value = "54683"
$ws.Range("B11").Value = "TON"
$ws.Range("C11").Value = "toncoin"
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = "-0.78"
$ws.Range("F11").Value = "9.50"
$ws.Range("G11").Value = "0.000092"
$ws.Range("H11").Value = 179048558.5134867
$ws.Range("I11").Value = "3468312277.00"
$ws.Range("J11").Value = "5047558528"

# Row 12
$ws.Range("A12").Value = "257"
$ws.Range("B12").Value = "ADA"
$ws.Range("C12").Value = "cardano"
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = "-2.01"
$ws.Range("F12").Value = "-0.59"
$ws.Range("G12").Value = "0.000007"
$ws.Range("H12").Value = 293771116.1444727
$ws.Range("I12").Value = "35489244418.00"
$ws.Range("J12").Value = "45000000000"

# Row 13
$ws.Range("A13").Value = "45088"
$ws.Range("B13").Value = "SHIB"
$ws.Range("C13").Value = "shiba-inu"
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = "-1.54"
$ws.Range("F13").Value = "-0.81"
$ws.Range("G13").Value = "3.74E-10"
$ws.Range("H13").Value = 65258366.33592062
$ws.Range("I13").Value = "589289410812691.00"
$ws.Range("J13").Value = ""

# Row 14
$ws.Range("A14").Value = "44883"
$ws.Range("B14").Value = "AVAX"
$ws.Range("C14").Value = "avalanche"
$ws.Range("D14").Value = 13
$ws.Range("E14").Value = "-1.10"
$ws.Range("F14").Value = "6.22"
$ws.Range("G14").Value = "0.000580"
$ws.Range("H14").Value = 408026138.0630771
$ws.Range("I14").Value = "377285092.00"
$ws.Range("J14").Value = "720000000"

# Row 15
$ws.Range("A15").Value = "2713"
$ws.Range("B15").Value = "TRX"
$ws.Range("C15").Value = "tron"
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = "1.55"
$ws.Range("F15").Value = "0.75"
$ws.Range("G15").Value = "0.000002"
$ws.Range("H15").Value = 356722804.212838
$ws.Range("I15").Value = "87923847381.00"
$ws.Range("J15").Value = ""

# Row 16
$ws.Range("A16").Value = "33422"
$ws.Range("B16").Value = "WBTC"
$ws.Range("C16").Value = "wrapped-bitcoin"
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = "0.33"
$ws.Range("F16").Value = "0.03"
$ws.Range("G16").Value = "0.998031"
$ws.Range("H16").Value = 72781170.93854488
$ws.Range("I16").Value = "155986.00"
$ws.Range("J16").Value = ""

# Row 17
$ws.Range("A17").Value = "2321"
$ws.Range("B17").Value = "BCH"
$ws.Range("C17").Value = "bitcoin-cash"
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = "1.37"
$ws.Range("F17").Value = "5.22"
$ws.Range("G17").Value = "0.007661"
$ws.Range("H17").Value = 373318882.4025794
$ws.Range("I17").Value = "19696959.00"
$ws.Range("J17").Value = "21000000"

# Row 18
$ws.Range("A18").Value = "45219"
$ws.Range("B18").Value = "DOT"
$ws.Range("C18").Value = "polkadot"
$ws.Range("D18").Value = 17
$ws.Range("E18").Value = "-0.25"
$ws.Range("F18").Value = "9.40"
$ws.Range("G18").Value = "0.000113"
$ws.Range("H18").Value = 216505614.6039383
$ws.Range("I18").Value = "1307052068.00"
$ws.Range("J18").Value = "1388001203.0802"

# Row 19
$ws.Range("A19").Value = "2751"
$ws.Range("B19").Value = "LINK"
$ws.Range("C19").Value = "chainlink"
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = "-1.79"
$ws.Range("F19").Value = "3.11"
$ws.Range("G19").Value = "0.000226"
$ws.Range("H19").Value = 279710790.4686785
$ws.Range("I19").Value = "587099970.00"
$ws.Range("J19").Value = "1000000000"

# Row 20
$ws.Range("A20").Value = "48563"
$ws.Range("B20").Value = "NEAR"
$ws.Range("C20").Value = "near-protocol"
$ws.Range("D20").Value = 19
$ws.Range("E20").Value = "2.06"
$ws.Range("F20").Value = "8.79"
$ws.Range("G20").Value = "0.000116"
$ws.Range("H20").Value = 438773796.3666959
$ws.Range("I20").Value = "1043761976.00"
$ws.Range("J20").Value = "1043761976"

# Row 21
$ws.Range("A21").Value = "1"
$ws.Range("B21").Value = "LTC"
$ws.Range("C21").Value = "litecoin"
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = "0.78"
$ws.Range("F21").Value = "-1.04"
$ws.Range("G21").Value = "0.001296"
$ws.Range("H21").Value = 501731536.8528561
$ws.Range("I21").Value = "74461469.00"
$ws.Range("J21").Value = "84000000"

# Row 22
$ws.Range("A22").Value = "33833"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "unus-sed-leo"
$ws.Range("D22").Value = 21
$ws.Range("E22").Value = "0.77"
$ws.Range("F22").Value = "-0.52"
$ws.Range("G22").Value = "0.000091"
$ws.Range("H22").Value = 3494352.94133756
$ws.Range("I22").Value = "927132386.00"
$ws.Range("J22").Value = ""

# Row 23
$ws.Range("A23").Value = "47305"
$ws.Range("B23").Value = "UNI"
$ws.Range("C23").Value = "uniswap"
$ws.Range("D23").Value = 22
$ws.Range("E23").Value = "0.44"
$ws.Range("F23").Value = "-1.10"
$ws.Range("G23").Value = "0.000119"
$ws.Range("H23").Value = 96116606.67136998
$ws.Range("I23").Value = "598736139.00"
$ws.Range("J23").Value = "1000000000"

# Row 24
$ws.Range("A24").Value = "118"
$ws.Range("B24").Value = "ETC"
$ws.Range("C24").Value = "ethereum-classic"
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = "0.62"
$ws.Range("F24").Value = "1.18"
$ws.Range("G24").Value = "0.000435"
$ws.Range("H24").Value = 400749386.7501777
$ws.Range("I24").Value = "145903895.00"
$ws.Range("J24").Value = "210700000"

# Row 25
$ws.Range("A25").Value = "44863"
$ws.Range("B25").Value = "RNDR"
$ws.Range("C25").Value = "render-token"
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = "10.58"
$ws.Range("F25").Value = "38.86"
$ws.Range("G25").Value = "0.000169"
$ws.Range("H25").Value = 309027174.4896991
$ws.Range("I25").Value = "374355803.00"
$ws.Range("J25").Value = "536870912"

# Row 26
$ws.Range("A26").Value = "33830"
$ws.Range("B26").Value = "ATOM"
$ws.Range("C26").Value = "cosmos"
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = "-0.90"
$ws.Range("F26").Value = "14.27"
$ws.Range("G26").Value = "0.000147"
$ws.Range("H26").Value = 154170773.6204464
$ws.Range("I26").Value = "389254388.00"
$ws.Range("J26").Value = ""

# Row 27
$ws.Range("A27").Value = "93841"
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "pepe"
$ws.Range("D27").Value = 26
$ws.Range("E27").Value = "-3.17"
$ws.Range("F27").Value = "12.32"
$ws.Range("G27").Value = "1.30E-10"
$ws.Range("H27").Value = 9423581.166320723
$ws.Range("I27").Value = "420689899999990.00"
$ws.Range("J27").Value = "420690000000000"

# Row 28
$ws.Range("A28").Value = "121595"
$ws.Range("B28").Value = "MNT"
$ws.Range("C28").Value = "mantle"
$ws.Range("D28").Value = 27
$ws.Range("E28").Value = "1.12"
$ws.Range("F28").Value = "3.38"
$ws.Range("G28").Value = "0.000017"
$ws.Range("H28").Value = 46941195.31005198
$ws.Range("I28").Value = "3231662126.00"
$ws.Range("J28").Value = "6219316795"

# Row 29
$ws.Range("A29").Value = "121613"
$ws.Range("B29").Value = "WIF"
$ws.Range("C29").Value = "dogwifhat"
$ws.Range("D29").Value = 28
$ws.Range("E29").Value = "-0.18"
$ws.Range("F29").Value = "27.64"
$ws.Range("G29").Value = "0.000053"
$ws.Range("H29").Value = 186818305.9748314
$ws.Range("I29").Value = "998920172.00"

# Row 30
$ws.Range("A30").Value = "111341"
$ws.Range("B30").Value = "APT"
$ws.Range("C30").Value = "aptos"
$ws.Range("D30").Value = 29
$ws.Range("E30").Value = "-0.50"
$ws.Range("F30").Value = "1.04"
$ws.Range("G30").Value = "0.000142"
$ws.Range("H30").Value = 118341616.0648232
$ws.Range("I30").Value = "368468672.00"
$ws.Range("J30").Value = "1084577363.9802"

# Row 31
$ws.Range("A31").Value = "48569"
$ws.Range("B31").Value = "STX"
$ws.Range("C31").Value = "stacks"
$ws.Range("D31").Value = 30
$ws.Range("E31").Value = "1.60"
$ws.Range("F31").Value = "-6.92"
$ws.Range("G31").Value = "0.000035"
$ws.Range("H31").Value = 102958474.3427878
$ws.Range("I31").Value = "1444838084.00"
$ws.Range("J31").Value = "1818000000"

# Row 32
$ws.Range("A32").Value = "89"
$ws.Range("B32").Value = "XLM"
$ws.Range("C32").Value = "stellar"
$ws.Range("D32").Value = 31
$ws.Range("E32").Value = "-1.98"
$ws.Range("F32").Value = "-1.62"
$ws.Range("G32").Value = "0.000002"
$ws.Range("H32").Value = 74872961.83501349
$ws.Range("I32").Value = "28919327940.00"
$ws.Range("J32").Value = "104303927518"

# Row 33
$ws.Range("A33").Value = "32607"
$ws.Range("B33").Value = "FIL"
$ws.Range("C33").Value = "filecoin"
$ws.Range("D33").Value = 32
$ws.Range("E33").Value = "-2.20"
$ws.Range("F33").Value = "1.87"
$ws.Range("G33").Value = "0.000094"
$ws.Range("H33").Value = 187135716.1690417
$ws.Range("I33").Value = "519800264.00"
$ws.Range("J33").Value = ""

# Row 34
$ws.Range("A34").Value = "100423"
$ws.Range("B34").Value = "FDUSD"
$ws.Range("C34").Value = "first-digital-usd"
$ws.Range("D34").Value = 33
$ws.Range("E34").Value = "0.04"
$ws.Range("F34").Value = "0.05"
$ws.Range("G34").Value = "0.000016"
$ws.Range("H34").Value = 5976846372.116722
$ws.Range("I34").Value = "3098764893.00"
